# Correct the "input/output overview" sheet to use the SF version values
# instead of the TPI version. This re-populates rows 2-11 (columns A-H)
# with the corrected, reordered data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, A (Input id), B (Value TOW), C, D, E, F, G (label), H (Smallest difference)
$data = @(
    @(2,  56100443116, 4988.32,  4988.3099999999995,  4988.83,              4988.2300000000005, 4988.889999999999,  "tax_calc_no_glob",  0.010000000000218279),
    @(3,  80080764310, 82.01,    110.21,               90.58,                90.58,               90.58,               "tax_calc_glob_p1",  8.569999999999993),
    @(4,  70062422978, 5457.95,  -6106.37,             -6106.37,             -6106.37,            -6106.37,            "tax_calc_no_glob",  11564.32),
    @(5,  91042947017, 71698.31, 71698.32,             71698.32,             72028.86,            71698.32,            "tax_calc_no_glob",  0.010000000009313226),
    @(6,  58101351196, 559.04,   1072.96,              -28294.34,            1074.3500000000004, -28294.34,            "tax_calc_no_glob",  513.9200000000001),
    @(7,  89051637777, 3627.39,  -4275.82,             -4275.82,             -3958.03,            -4275.82,            "tax_calc_glob_p2",  7585.42),
    @(8,  68042050362, -72081.72,-72049.04000000001,   -72049.04000000001,   -72046.83,           -72049.04000000001,  "tax_calc_no_glob",  32.679999999993015),
    @(9,  59093055145, 3591.83,  1535.6500000000005,   -41497.71000000001,   2411.2799999999997,  -52406.91,           "tax_calc_glob_p2",  1180.5500000000002),
    @(10, 79122109350, -2038.91, -2042.75,             -2003.41,             -1724.96,            -2003.41,            "tax_calc_no_glob",  3.839999999999918),
    @(11, 80091134994, 957.75,   -2050.44,             -2050.44,             -1732.65,            -2050.44,            "tax_calc_glob_p2",  2690.4)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
